# Rearrange the Application.xlsx HR sheet columns/header order and refresh
# data accordingly (PassCv/Interview/Hire block moved before the applicant
# detail columns), matching the re-imported Power BI layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 'ApplicationId'
$ws.Cells.Item(1, 2).Value = 'ApplicantId'
$ws.Cells.Item(1, 3).Value = 'FullName'
$ws.Cells.Item(1, 4).Value = 'DateOfBirth'
$ws.Cells.Item(1, 5).Value = 'Gender'
$ws.Cells.Item(1, 6).Value = 'Address'
$ws.Cells.Item(1, 7).Value = 'PhoneNumber'
$ws.Cells.Item(1, 8).Value = 'Email'
$ws.Cells.Item(1, 9).Value = 'Level'
$ws.Cells.Item(1, 10).Value = 'RecruitmentChannelId'
$ws.Cells.Item(1, 11).Value = 'JobPositionId'
$ws.Cells.Item(1, 12).Value = 'ApplicationDate'
$ws.Cells.Item(1, 13).Value = 'PassCv'
$ws.Cells.Item(1, 14).Value = 'InterviewDate'
$ws.Cells.Item(1, 15).Value = 'InterviewResult'
$ws.Cells.Item(1, 16).Value = 'HireDate'

# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 'Nguyễn Văn A'
$ws.Cells.Item(2, 4).Value = 32918
$ws.Cells.Item(2, 5).Value = 'Male'
$ws.Cells.Item(2, 6).Value = '123 Đường Láng, Hà Nội'
$ws.Cells.Item(2, 7).Value = 987654321
$ws.Cells.Item(2, 8).Value = 'nguyen.a@example.com'
$ws.Cells.Item(2, 9).Value = 'Senior'
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 45296
$ws.Cells.Item(2, 13).Value = 'Pass'
$ws.Cells.Item(2, 14).Value = 45301
$ws.Cells.Item(2, 15).Value = 'Pass'
$ws.Cells.Item(2, 16).Value = 45323

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = 'Trần Thị B'
$ws.Cells.Item(3, 4).Value = 33746
$ws.Cells.Item(3, 5).Value = 'Female'
$ws.Cells.Item(3, 6).Value = '456 Nguyễn Trãi, Hà Nội'
$ws.Cells.Item(3, 7).Value = 912345678
$ws.Cells.Item(3, 8).Value = 'tran.b@example.com'
$ws.Cells.Item(3, 9).Value = 'Junior'
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 45334
$ws.Cells.Item(3, 13).Value = 'Pass'
$ws.Cells.Item(3, 14).Value = 45340
$ws.Cells.Item(3, 15).Value = 'Pass'
$ws.Cells.Item(3, 16).Value = 45352

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 'Lê Văn C'
$ws.Cells.Item(4, 4).Value = 32450
$ws.Cells.Item(4, 5).Value = 'Male'
$ws.Cells.Item(4, 6).Value = '789 Cầu Giấy, Hà Nội'
$ws.Cells.Item(4, 7).Value = 901234567
$ws.Cells.Item(4, 8).Value = 'le.c@example.com'
$ws.Cells.Item(4, 9).Value = 'Fresher'
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 45376
$ws.Cells.Item(4, 13).Value = 'Fail'
$ws.Cells.Item(4, 14).Value = 45381
$ws.Cells.Item(4, 15).Value = 'Fail'
$ws.Cells.Item(4, 16).Value = 45352

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(5, 3).Value = 'Phạm Minh D'
$ws.Cells.Item(5, 4).Value = 34897
$ws.Cells.Item(5, 5).Value = 'Male'
$ws.Cells.Item(5, 6).Value = '101 Đường Trần Duy Hưng'
$ws.Cells.Item(5, 7).Value = 981122334
$ws.Cells.Item(5, 8).Value = 'pham.d@example.com'
$ws.Cells.Item(5, 9).Value = 'Intern'
$ws.Cells.Item(5, 10).Value = 4
$ws.Cells.Item(5, 11).Value = 4
$ws.Cells.Item(5, 12).Value = 45384
$ws.Cells.Item(5, 13).Value = 'Pass'
$ws.Cells.Item(5, 14).Value = 45390
$ws.Cells.Item(5, 15).Value = 'Pass'
$ws.Cells.Item(5, 16).Value = 45413

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 5
$ws.Cells.Item(6, 3).Value = 'Vũ Thị E'
$ws.Cells.Item(6, 4).Value = 33491
$ws.Cells.Item(6, 5).Value = 'Female'
$ws.Cells.Item(6, 6).Value = '202 Láng Hạ, Hà Nội'
$ws.Cells.Item(6, 7).Value = 976543210
$ws.Cells.Item(6, 8).Value = 'vu.e@example.com'
$ws.Cells.Item(6, 9).Value = 'Senior'
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 5
$ws.Cells.Item(6, 12).Value = 45428
$ws.Cells.Item(6, 13).Value = 'Pass'
$ws.Cells.Item(6, 14).Value = 45433
$ws.Cells.Item(6, 15).Value = 'Pass'
$ws.Cells.Item(6, 16).Value = 45444

# Row 7
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 6
$ws.Cells.Item(7, 3).Value = 'Đỗ Văn F'
$ws.Cells.Item(7, 4).Value = 32132
$ws.Cells.Item(7, 5).Value = 'Male'
$ws.Cells.Item(7, 6).Value = '303 Nguyễn Xiển, Hà Nội'
$ws.Cells.Item(7, 7).Value = 965432109
$ws.Cells.Item(7, 8).Value = 'do.f@example.com'
$ws.Cells.Item(7, 9).Value = 'Junior'
$ws.Cells.Item(7, 10).Value = 2
$ws.Cells.Item(7, 11).Value = 6
$ws.Cells.Item(7, 12).Value = 45461
$ws.Cells.Item(7, 13).Value = 'Fail'
$ws.Cells.Item(7, 14).Value = 45468
$ws.Cells.Item(7, 15).Value = 'Cancel'
$ws.Cells.Item(7, 16).Value = 45352

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 7
$ws.Cells.Item(8, 3).Value = 'Nguyễn Thị G'
$ws.Cells.Item(8, 4).Value = 34429
$ws.Cells.Item(8, 5).Value = 'Female'
$ws.Cells.Item(8, 6).Value = '404 Hoàng Quốc Việt'
$ws.Cells.Item(8, 7).Value = 954321098
$ws.Cells.Item(8, 8).Value = 'nguyen.g@example.com'
$ws.Cells.Item(8, 9).Value = 'Intern'
$ws.Cells.Item(8, 10).Value = 3
$ws.Cells.Item(8, 11).Value = 7
$ws.Cells.Item(8, 12).Value = 45484
$ws.Cells.Item(8, 13).Value = 'Pass'
$ws.Cells.Item(8, 14).Value = 45488
$ws.Cells.Item(8, 15).Value = 'Pass'
$ws.Cells.Item(8, 16).Value = 45505

# Row 9
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 8
$ws.Cells.Item(9, 3).Value = 'Hoàng Văn H'
$ws.Cells.Item(9, 4).Value = 33114
$ws.Cells.Item(9, 5).Value = 'Male'
$ws.Cells.Item(9, 6).Value = '505 Kim Mã, Hà Nội'
$ws.Cells.Item(9, 7).Value = 943210987
$ws.Cells.Item(9, 8).Value = 'hoang.h@example.com'
$ws.Cells.Item(9, 9).Value = 'Fresher'
$ws.Cells.Item(9, 10).Value = 4
$ws.Cells.Item(9, 11).Value = 8
$ws.Cells.Item(9, 12).Value = 45507
$ws.Cells.Item(9, 13).Value = 'Pass'
$ws.Cells.Item(9, 14).Value = 45514
$ws.Cells.Item(9, 15).Value = 'Pass'
$ws.Cells.Item(9, 16).Value = 45536

# Row 10
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 9
$ws.Cells.Item(10, 3).Value = 'Lý Thị I'
$ws.Cells.Item(10, 4).Value = 34047
$ws.Cells.Item(10, 5).Value = 'Female'
$ws.Cells.Item(10, 6).Value = '606 Đội Cấn, Hà Nội'
$ws.Cells.Item(10, 7).Value = 932109876
$ws.Cells.Item(10, 8).Value = 'ly.i@example.com'
$ws.Cells.Item(10, 9).Value = 'Senior'
$ws.Cells.Item(10, 10).Value = 1
$ws.Cells.Item(10, 11).Value = 9
$ws.Cells.Item(10, 12).Value = 45542
$ws.Cells.Item(10, 13).Value = 'Pass'
$ws.Cells.Item(10, 14).Value = 45550
$ws.Cells.Item(10, 15).Value = 'Fail'
$ws.Cells.Item(10, 16).Value = 45352

# Row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 10
$ws.Cells.Item(11, 3).Value = 'Bùi Văn K'
$ws.Cells.Item(11, 4).Value = 35246
$ws.Cells.Item(11, 5).Value = 'Male'
$ws.Cells.Item(11, 6).Value = '707 Phạm Hùng, Hà Nội'
$ws.Cells.Item(11, 7).Value = 921098765
$ws.Cells.Item(11, 8).Value = 'bui.k@example.com'
$ws.Cells.Item(11, 9).Value = 'Junior'
$ws.Cells.Item(11, 10).Value = 2
$ws.Cells.Item(11, 11).Value = 10
$ws.Cells.Item(11, 12).Value = 45587
$ws.Cells.Item(11, 13).Value = 'Pass'
$ws.Cells.Item(11, 14).Value = 45592
$ws.Cells.Item(11, 15).Value = 'Pass'
$ws.Cells.Item(11, 16).Value = 45611

# Update the active cell selection to match the saved view state
$ws.Range("L18").Select()
